# "Add files via upload" — a new leetcoder entry (rank 1762, Silvia42) is
# inserted into the leaderboard as the new row 8, pushing every following
# row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Insert a new blank row at position 8 (shifts rows 8..23 down to 9..24,
# inheriting number formatting from the row above just like Excel does).
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new entrant's stats.
$ws.Cells.Item(8, 1).Value = 1762
$ws.Cells.Item(8, 2).Value = "https://leetcode.com/u/Silvia42/"
$ws.Cells.Item(8, 3).Value = 56
$ws.Cells.Item(8, 4).Value = 1625
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 11
$ws.Cells.Item(8, 11).Value = 13
$ws.Cells.Item(8, 12).Value = "No data"

# Match the saved view state: selection on M8.
$ws.Range("M8").Select() | Out-Null
